$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffix = ".age_trait"

# Columns B through P on row 1 all get the ".age_trait" suffix appended
# to their existing header text (column A "Country" is left unchanged).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")

foreach ($col in $cols) {
    $cell = $ws.Range("$col`1")
    $current = $cell.Value2
    $cell.Value = "$current$suffix"
}
